# Weekly update: insert two new observation rows into the "Arveja Verde"
# (Terminal Hortofrutícola Agro Chillán) dataset, pushing the existing rows
# below them down, matching the upstream weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 21 (existing rows 21.. shift to 22..) ---
$ws.Rows.Item(21).Insert()

$ws.Cells.Item(21, 1).Value = 7
$ws.Cells.Item(21, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21, 3).Value = "Ñuble"
$ws.Cells.Item(21, 4).Value = 44554
$ws.Cells.Item(21, 5).Value = 16
$ws.Cells.Item(21, 6).Value = 100112022
$ws.Cells.Item(21, 7).Value = "Arveja Verde"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 120
$ws.Cells.Item(21, 11).Value = 15000
$ws.Cells.Item(21, 12).Value = 16000
$ws.Cells.Item(21, 13).Value = 15500
$ws.Cells.Item(21, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(21, 16).Value = 620
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# --- Insert second new row at position 32 (existing rows 32.. shift to 33..) ---
$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(32, 3).Value = "Ñuble"
$ws.Cells.Item(32, 4).Value = 44553
$ws.Cells.Item(32, 5).Value = 16
$ws.Cells.Item(32, 6).Value = 100112022
$ws.Cells.Item(32, 7).Value = "Arveja Verde"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 14000
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = 14500
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 580
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Keep the date-formatted column (D) consistent with the rest of the table.
$ws.Range("D21").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("D32").NumberFormat = $ws.Range("D31").NumberFormat
